$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Quantum Revelations: Unraveling Nature's Enigmatic Secrets" "Delving into the World of Politics: A Comprehensive Guide"

# --- Author name ---
Replace-Text "Oliver Saunders" "Ethan Bennett"

# --- Email address (two runs keep their own identity: "oliver" -> "bennetthan@gnail";
#     then "saunders96@emailcentral" + "." + "net" collapse into "cam") ---
Replace-Text "oliver" "bennetthan@gnail"
Replace-Text "saunders96@emailcentral.net" "cam"

# The engine coalesces same-format runs in a touched paragraph; force a re-split so the
# middle "." stays its own run (matching the source run layout) by toggling Bold off/on
# (net no-op) on just that character.
$fr = $d.Content
$fr.Find.Execute("bennetthan@gnail", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dot = $d.Range($fr.End, $fr.End + 1)
$dot.Bold = 1
$dot.Bold = 0

# --- First body paragraph, block 1 (before the first double line-break) ---
Replace-Text "As we embark on a captivating odyssey into the realm of quantum mechanics, a thought-provoking dance of particles and waves, the very foundations of our universe are unveiled. The subatomic realm, an enigmatic tapestry woven with probability and uncertainty, beckons us to explore phenomena that defy classical intuition. From the intrinsic interconnectedness of entangled particles to the tunnel-like passage of particles through impassable barriers, quantum mechanics has revolutionized our understanding of the universe and profoundly influenced diverse fields, ranging from computing to cryptography." "Politics, a multifaceted and ever-evolving sphere of human interaction, holds immense significance in shaping our world. It encompasses the systems, institutions, and processes through which societies are governed, decisions are made, and resources are allocated. Understanding politics is crucial for young minds as it equips them with the knowledge and skills necessary to navigate the intricacies of a complex and interconnected world."

# --- First body paragraph, block 2 (between the two double line-breaks) ---
Replace-Text "In the quantum domain, the classical laws governing the macroscopic world unravel, leaving us with a fascinating enigma. Particles exhibit paradoxical behaviors, simultaneously existing and occupying multiple states, an ethereal symphony of possibilities. These mind-bending phenomena have led to the formulation of profound interpretations, questioning the fabric of reality itself. The tension between competing interpretations gives rise to lively debates and captivating paradoxes, challenging our conventional notions of time, locality, and determinism. The quantum realm beckons us to confront these paradoxical intricacies, pushing the boundaries of human knowledge and understanding." "In this essay, we will embark on a journey to explore the multifaceted nature of politics. We will delve into its fundamental principles, analyze different political ideologies, and examine the diverse structures and institutions that govern nations. By gaining a deeper understanding of politics, we can foster active and informed citizenship, enabling individuals to participate effectively in shaping their communities and societies."

# --- First body paragraph, block 3 (after the second double line-break) ---
Replace-Text "Quantum mechanics has propelled advancements across disciplines. In the realm of cryptography, it holds the promise of unbreakable encryption, as the inherent uncertainty of quantum systems renders eavesdropping impossible. This transformative potential has inspired the creation of quantum communication networks, poised to revolutionize secure information transfer. Furthermore, quantum computing, harnessing the power of quantum entanglement and superposition, offers the prospect of exponential gains in computational capacity, opening uncharted frontiers for simulations and optimizations. These transformative applications are poised to usher in an era of unprecedented technological progress and redefine the boundaries of human achievement." "The study of politics begins with an exploration of its core concepts. We will investigate the nature of power, authority, and legitimacy, delving into the various ways in which individuals and institutions exercise influence over others. We will examine the relationship between politics and economics, considering how economic structures and policies impact political decision-making. Additionally, we will analyze the role of political culture and ideology in shaping the political landscape of societies."

# --- "Summary" heading: self-replace to drop the stray lastRenderedPageBreak marker ---
Replace-Text "Summary" "Summary"

# --- Summary paragraph ---
Replace-Text "In this exploration of quantum mechanics, we have delved into the perplexing realm of subatomic particles, uncovering phenomena that challenge conventional wisdom. Quantum mechanics has redefined our understanding of nature, fostering lively debates about reality's fundamental structure. It has fueled remarkable innovations in computing, cryptography, and other fields, hinting at a future brimming with possibilities. As we continue to unravel the enigmatic secrets of the quantum realm, we may one day resolve its paradoxes and forge a deeper connection with the universe's enigmatic tapestry." "This essay provides a comprehensive overview of politics, exploring its fundamental principles, ideologies, and structures. It highlights the importance of political engagement and emphasizes the role of individuals in shaping political outcomes. By understanding the complexities of politics, young minds can become active and informed citizens, capable of making informed decisions and contributing positively to their communities and societies."

# --- Add a new trailing empty paragraph before the section break ---
$d.Content.InsertParagraphAfter()
